$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 130886843
$ws.Range("Q2").Value = 434321
$ws.Range("R2").Value = 7052458

# Row 3
$ws.Range("A3").Value = 130886842
$ws.Range("Q3").Value = 434316
$ws.Range("R3").Value = 7052462

# Row 13
$ws.Range("A13").Value = 130886774
$ws.Range("B13").Value = 57884
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = "Tretåig hackspett"
$ws.Range("G13").Value = "Picoides tridactylus"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("Q13").Value = 434519
$ws.Range("R13").Value = 7052220
$ws.Range("AC13").Value = "Ringhack"

# Row 14
$ws.Range("A14").Value = 130886839
$ws.Range("B14").Value = 91828
$ws.Range("E14").Value = 5432
$ws.Range("F14").Value = "Granticka"
$ws.Range("G14").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H14").Value = ""
$ws.Range("Q14").Value = 434382
$ws.Range("R14").Value = 7052442
$ws.Range("AC14").ClearContents()

# Row 21
$ws.Range("A21").Value = 130886831
$ws.Range("B21").Value = 57881
$ws.Range("E21").Value = 100049
$ws.Range("F21").Value = "Spillkråka"
$ws.Range("G21").Value = "Dryocopus martius"
$ws.Range("H21").Value = "(Linnaeus, 1758)"
$ws.Range("Q21").Value = 434871
$ws.Range("R21").Value = 7051709
$ws.Range("AC21").Value = "Hack"

# Row 22
$ws.Range("A22").Value = 130886836
$ws.Range("B22").Value = 91828
$ws.Range("E22").Value = 5432
$ws.Range("F22").Value = "Granticka"
$ws.Range("G22").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H22").Value = ""
$ws.Range("Q22").Value = 434854
$ws.Range("R22").Value = 7051718
$ws.Range("AC22").ClearContents()

# Row 23
$ws.Range("A23").Value = 130886792
$ws.Range("Q23").Value = 434158
$ws.Range("R23").Value = 7052168
$ws.Range("AC23").Value = "Ringhack äldre"

# Row 24
$ws.Range("A24").Value = 130886801
$ws.Range("Q24").Value = 434001
$ws.Range("R24").Value = 7052192
$ws.Range("AC24").Value = "Ringhack färska och äldre"

# Row 25
$ws.Range("A25").Value = 130886793
$ws.Range("Q25").Value = 434143
$ws.Range("R25").Value = 7052197
$ws.Range("AC25").Value = "Ringhack äldre"

# Row 26
$ws.Range("A26").Value = 130886794
$ws.Range("Q26").Value = 434140
$ws.Range("R26").Value = 7052192
$ws.Range("AC26").Value = "Ringhack färska och äldre"

# Row 27
$ws.Range("A27").Value = 130886818
$ws.Range("B27").Value = 57884
$ws.Range("E27").Value = 100109
$ws.Range("F27").Value = "Tretåig hackspett"
$ws.Range("G27").Value = "Picoides tridactylus"
$ws.Range("Q27").Value = 434272
$ws.Range("R27").Value = 7052031
$ws.Range("AC27").Value = "Ringhack äldre"

# Row 39
$ws.Range("A39").Value = 130886825
$ws.Range("B39").Value = 57884
$ws.Range("E39").Value = 100109
$ws.Range("F39").Value = "Tretåig hackspett"
$ws.Range("G39").Value = "Picoides tridactylus"
$ws.Range("H39").Value = "(Linnaeus, 1758)"
$ws.Range("Q39").Value = 434476
$ws.Range("R39").Value = 7051885
$ws.Range("AC39").Value = "Ringhack äldre"

# Row 40
$ws.Range("A40").Value = 130886805
$ws.Range("Q40").Value = 433991
$ws.Range("R40").Value = 7052188
$ws.Range("AC40").Value = "Ringhack"

# Row 41
$ws.Range("A41").Value = 130886845
$ws.Range("B41").Value = 91828
$ws.Range("E41").Value = 5432
$ws.Range("F41").Value = "Granticka"
$ws.Range("G41").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H41").Value = ""
$ws.Range("Q41").Value = 434162
$ws.Range("R41").Value = 7052165
$ws.Range("AC41").ClearContents()

# Row 51
$ws.Range("A51").Value = 130886837
$ws.Range("B51").Value = 91828
$ws.Range("E51").Value = 5432
$ws.Range("F51").Value = "Granticka"
$ws.Range("G51").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H51").Value = ""
$ws.Range("Q51").Value = 434513
$ws.Range("R51").Value = 7052004
$ws.Range("AC51").ClearContents()

# Row 52
$ws.Range("A52").Value = 130886791
$ws.Range("Q52").Value = 434162
$ws.Range("R52").Value = 7052153

# Row 53
$ws.Range("A53").Value = 130886765
$ws.Range("Q53").Value = 434505
$ws.Range("R53").Value = 7052005
$ws.Range("AC53").Value = "Ringhack"

# Row 54
$ws.Range("A54").Value = 130886788
$ws.Range("Q54").Value = 434171
$ws.Range("R54").Value = 7052213

# Row 55
$ws.Range("A55").Value = 130886826
$ws.Range("Q55").Value = 434489
$ws.Range("R55").Value = 7051863
$ws.Range("AC55").Value = "Ringhack äldre"

# Row 56
$ws.Range("A56").Value = 130886811
$ws.Range("Q56").Value = 434077
$ws.Range("R56").Value = 7052133
$ws.Range("AC56").Value = "Ringhack färska"

# Row 57
$ws.Range("A57").Value = 130886785
$ws.Range("B57").Value = 57884
$ws.Range("E57").Value = 100109
$ws.Range("F57").Value = "Tretåig hackspett"
$ws.Range("G57").Value = "Picoides tridactylus"
$ws.Range("H57").Value = "(Linnaeus, 1758)"
$ws.Range("Q57").Value = 434191
$ws.Range("R57").Value = 7052193
$ws.Range("AC57").Value = "Ringhack äldre"

# Row 73
$ws.Range("A73").Value = 130886823
$ws.Range("B73").Value = 57884
$ws.Range("D73").Value = "NT"
$ws.Range("E73").Value = 100109
$ws.Range("F73").Value = "Tretåig hackspett"
$ws.Range("G73").Value = "Picoides tridactylus"
$ws.Range("I73").Value = ""
$ws.Range("Q73").Value = 434499
$ws.Range("R73").Value = 7051916
$ws.Range("AC73").Value = "Ringhack"
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("K73").ClearContents()
$ws.Range("L73").ClearContents()

# Row 74
$ws.Range("A74").Value = 130886813
$ws.Range("Q74").Value = 434112
$ws.Range("R74").Value = 7052117

# Row 75
$ws.Range("A75").Value = 130886762
$ws.Range("Q75").Value = 434867
$ws.Range("R75").Value = 7051762
$ws.Range("AC75").Value = "Ringhack äldre"

# Row 76
$ws.Range("A76").Value = 130886821
$ws.Range("Q76").Value = 434468
$ws.Range("R76").Value = 7051906

# Row 77
$ws.Range("A77").Value = 130886789
$ws.Range("Q77").Value = 434159
$ws.Range("R77").Value = 7052197

# Row 78
$ws.Range("A78").Value = 130886832
$ws.Range("B78").Value = 57988
$ws.Range("D78").Value = "LC"
$ws.Range("E78").Value = 103031
$ws.Range("F78").Value = "Lavskrika"
$ws.Range("G78").Value = "Perisoreus infaustus"
$ws.Range("I78").Value = "1"
$ws.Range("K78").Value = ""
$ws.Range("L78").Value = ""
$ws.Range("M78").Value = "födosökande"
$ws.Range("N78").Value = "observerad"
$ws.Range("Q78").Value = 434123
$ws.Range("R78").Value = 7052111
$ws.Range("AC78").ClearContents()
